{"js": "// Commit: \"Remove semicolon in equation area\"\n// The variable-definition (\"equation area\") lines in this document each read\n// like \"<name> = <description>; <value> ...\" \u2014 every semicolon is swapped\n// for a colon (the space that already follows each semicolon is left as-is,\n// so \"; \" naturally becomes \": \"):\n//   Size of air compressor; ${HP}                           -> ...: ${HP}\n//   Compressor power fraction; ${FR}                         -> ...: ${FR}\n//   Fraction of electrical energy converted into heat; ${EC} -> ...: ${EC}\n//   Conversion factor from HP to MMBtu/hr; 0.002544          -> ...: 0.002544\n//   Efficiency of heat recovery; estimated ${EHR}            -> ...: estimated ${EHR}\n//   Wintertime operating hours; ${OH}                        -> ...: ${OH}\n//\n// Every semicolon in the body belongs to one of these six lines (verified\n// against the source document), so a body-wide search/replace of the bare\n// \";\" character is sufficient and unambiguous.\n\nconst body = context.document.body;\nconst results = body.search(\";\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\":\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Commit: \"Remove semicolon in equation area\"\n#\n# The variable-definition (\"equation area\") lines in this document each read\n# like \"<name> = <description>; <value> ...\" - every semicolon is swapped for\n# a colon (the space that already follows each semicolon is left as-is, so\n# \"; \" naturally becomes \": \"):\n#   Size of air compressor; ${HP}                           -> ...: ${HP}\n#   Compressor power fraction; ${FR}                         -> ...: ${FR}\n#   Fraction of electrical energy converted into heat; ${EC} -> ...: ${EC}\n#   Conversion factor from HP to MMBtu/hr; 0.002544          -> ...: 0.002544\n#   Efficiency of heat recovery; estimated ${EHR}            -> ...: estimated ${EHR}\n#   Wintertime operating hours; ${OH}                        -> ...: ${OH}\n#\n# Every semicolon in the body belongs to one of these six lines (verified\n# against the source document), so a document-wide Find/Replace of the bare\n# \";\" character is sufficient and unambiguous.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Execute(\n  \";\",    # FindText\n  $false, # MatchCase\n  $false, # MatchWholeWord\n  $false, # MatchWildcards\n  $false, # MatchSoundsLike\n  $false, # MatchAllWordForms\n  $true,  # Forward\n  1,      # Wrap            (wdFindContinue)\n  $false, # Format\n  \":\",    # ReplaceWith\n  2       # Replace         (wdReplaceAll)\n) | Out-Null\n"}
